# ACLXRAY Lab deploy instructions.docx - apply commit "set time zone and
# delete strorage account"
#
# Content-level changes applied:
#  1. "Shutdown Time Zone (default 'Eastern Standard Time')" bullet is
#     reworded to "Shutdown Time Zone - default is a current time zone of
#     the user" (the hard-coded default time zone wording is replaced
#     because the script no longer hard codes Eastern Standard Time).
#  2. The matching "Default Value: 'Eastern Standard Time'" bullet under the
#     -shutdownTimeZone parameter reference section is removed completely,
#     since that parameter no longer has a fixed default value.
#
# (The rest of the diff is Word's own grammar-checker <w:proofErr/> markers
# being cleared out and the runs they used to split merging back together;
# that is a side effect of Word's proofing pass and carries no visible text
# change, so there is nothing further to type/replace for it.)

$d = $word.ActiveDocument

# Non-breaking space, matches the character actually used throughout this
# document's body text (pasted-from-web content uses U+00A0 instead of a
# plain space in many runs).
$nbsp = [char]0x00A0
$endash = [char]0x2013

# --- 1. Reword the "Shutdown Time Zone" bullet -----------------------------
$oldTimeZoneBullet = "Shutdown${nbsp}Time${nbsp}Zone${nbsp}(default${nbsp}'Eastern${nbsp}Standard${nbsp}Time')"
$newTimeZoneBullet = "Shutdown${nbsp}Time${nbsp}Zone $endash default is a current time zone of the user"

$found = $d.Content.Find.Execute($oldTimeZoneBullet, $true, $false, $false, $false, $false, $true, 1, $false, $newTimeZoneBullet, 2)
if (-not $found) {
    throw "Could not find the 'Shutdown Time Zone' bullet to reword"
}

# --- 2. Delete the "Default Value: 'Eastern Standard Time'" bullet ---------
$oldDefaultValueBullet = "Default${nbsp}Value:${nbsp}'Eastern${nbsp}Standard${nbsp}Time'"

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $oldDefaultValueBullet) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Default Value: Eastern Standard Time' paragraph to delete"
}

$d.Paragraphs.Item($targetIndex).Range.Delete()
